# Fill in blanks in the HID code map (Sheet2, which is the worksheet
# holding the "keypresses" table -- note the tab literally named
# "Sheet1" in this workbook is the other, empty sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New HID usage notes for previously-blank cells.
$ws.Range("D65").Value = "HID_KEYBOARD_MENU"
$ws.Range("D67").Value = "HID_KEYBOARD_SELECT"
$ws.Range("D69").Value = "HID_KEYBOARD_EXECUTE"

$ws.Range("D115").Value = "HID_KEYBOARD_MUTE"

$ws.Range("D116").Value = "HID_KEYBOARD_VOLUME_DOWN"
$ws.Range("E116").Value = "BRIGHTNESS?"

$ws.Range("D117").Value = "HID_KEYBOARD_VOLUME_UP"
$ws.Range("E117").Value = "BRIGHTNESS?"

$ws.Range("D118").Value = "HID_KEYBOARD_POWER"
$ws.Range("E118").Value = "SHOULD THIS BE A KEY??"

# The table was originally backed by a Power Query ("keypresses"); break
# that link so it becomes a plain table (drops queryTable-only attributes
# such as tableType/queryTableFieldId/uniqueName) and removes the
# associated hidden ExternalData_1 defined name.
$lo = $ws.ListObjects.Item(1)
$lo.Unlink()

foreach ($n in $wb.Names) {
    if ($n.Name -like "*ExternalData_1*") {
        $n.Delete()
    }
}

# Match the final selection/scroll position left behind in the sheet.
$ws.Range("D118").Select() | Out-Null
